$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete lichtwark values in B2, D2, E2; update C2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.3319794989134781
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 tweaks
$ws.Range("B3").Value = 5.6375100864256718
$ws.Range("C3").Value = 7.0164431192135419
$ws.Range("D3").Value = 8.7406576949142938
$ws.Range("E3").Value = 4.0263754227963036

# Update the selected range to match the new area of interest
$ws.Range("B1:E3").Select()
